$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header-row captions so each column is suffixed with the
#    respective format-version name instead of the generic "_old"/"_new".
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# column K (11) stays "diff" - untouched

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# 2. Freeze the header row (split/freeze after row 1, top-left of the
#    scrollable area becomes A2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into a real Excel Table ("Table1") with banded
#    rows and an AutoFilter on the header row.
$usedRange = $ws.Range("A1:U88")
$table = $ws.ListObjects.Add(1, $usedRange, [System.Type]::Missing, 1)
$table.Name = "Table1"
